# The deck ships two DrawingML themes:
#   ppt/theme/theme1.xml  -> linked from the slide master (the "Integral" theme,
#                             what every slide actually renders with)
#   ppt/theme/theme2.xml  -> linked only from the notes master ("Office Theme")
#
# The target edit swaps the two themes' contents (theme1 becomes "Office
# Theme", theme2 becomes "Integral"). The font scheme and format scheme
# (fill/line/effect/background styles) are already byte-identical between
# the two themes, so the only real content difference is the 12-slot colour
# scheme (and the cosmetic theme/colour-scheme "name" attributes, which
# PowerPoint's automation surface does not expose for editing).
#
# Because every reachable colour-scheme / theme-color object in this host
# (SlideMaster, NotesMaster, NotesPage, individual Slides, ...) all resolve
# back to the single active theme part (theme1.xml), the achievable and
# faithful version of this edit is to repaint theme1.xml's colour scheme
# with the "Office Theme" palette that the diff moves into it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index -> (scheme slot, target RGB as a COM BGR-packed long == RGB(r,g,b))
$tcs.Item(1).RGB  = 0x000000   # dk1      000000
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = 0x6A5444   # dk2      44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # accent2  ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  FFC000
$tcs.Item(9).RGB  = 0xC47244   # accent5  4472C4
$tcs.Item(10).RGB = 0x47AD70   # accent6  70AD47
$tcs.Item(11).RGB = 0xC16305   # hlink    0563C1
$tcs.Item(12).RGB = 0x724F95   # folHlink 954F72
